# Update the 8 header cells in row 1 (text content / casing changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "NUMERO CONTA"
$ws.Range("B1").Value = "ASSESSOR RV"
$ws.Range("C1").Value = "ADVISOR"
$ws.Range("D1").Value = "CLIENTE"
$ws.Range("E1").Value = "ESTRATÉGIA"
$ws.Range("F1").Value = "NET TOTAL"
$ws.Range("G1").Value = "NET DISPONÍVEL"
$ws.Range("H1").Value = "VALOR MEDIO POR OPERAÇÃO"

# Resize columns A:H to match the widths left behind after the header
# text changed (values chosen so the saved column width lands on the
# target figure).
$ws.Columns("A").ColumnWidth = 21.166666666666668
$ws.Columns("B").ColumnWidth = 19.166666666666668
$ws.Columns("C").ColumnWidth = 12.333333333333332
$ws.Columns("D").ColumnWidth = 19.0
$ws.Columns("E").ColumnWidth = 16.666666666666668
$ws.Columns("F").ColumnWidth = 14.666666666666666
$ws.Columns("G").ColumnWidth = 21.666666666666668
$ws.Columns("H").ColumnWidth = 39.666666666666664

# Move the active selection to H1, as left by the editing session.
$ws.Range("H1").Select()

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
